# scheduling.xlsx — first stage of adapting to project's problem
#
# - Rename "Available" -> "Instructors"
# - Insert a new "Students" sheet (between Instructors and Required) and
#   move the student-crew / next-event table from the Instructors sheet
#   (columns G/I) onto it, adding MinHours/MaxHours columns.
# - Un-hide the MinHours/MaxHours columns (D/E) on the Instructors sheet,
#   clear the now-relocated columns G/I, and bump a handful of required
#   staff-counts.
# - Populate the (until now empty) "Sheet3" with the Event /
#   ValueForCompleting lookup table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename "Available" -> "Instructors"
# ---------------------------------------------------------------------
$instructors = $wb.Worksheets.Item("Available")
$instructors.Name = "Instructors"

# ---------------------------------------------------------------------
# 2. Insert new "Students" sheet right after "Instructors"
# ---------------------------------------------------------------------
$students = $wb.Worksheets.Add($null, $instructors)
$students.Name = "Students"

$students.Cells.Item(1,1).Value = "StudentCrewName"
$students.Cells.Item(1,2).Value = "MinHours"
$students.Cells.Item(1,3).Value = "MaxHours"
$students.Cells.Item(1,4).Value = "NextEvent"

$students.Cells.Item(2,1).Value = "Chambers, Button"
$students.Cells.Item(2,2).Value = 1
$students.Cells.Item(2,3).Value = 1
$students.Cells.Item(2,4).Value = "CAS-3"

$students.Cells.Item(3,1).Value = "Solano, Sorensen"
$students.Cells.Item(3,2).Value = 1
$students.Cells.Item(3,3).Value = 1
$students.Cells.Item(3,4).Value = "SCAR-2"

$students.Cells.Item(4,1).Value = "Goins, Burrell"
$students.Cells.Item(4,2).Value = 1
$students.Cells.Item(4,3).Value = 1
$students.Cells.Item(4,4).Value = "TR-4"

$students.Cells.Item(5,1).Value = "Huggins, Cornett"
$students.Cells.Item(5,2).Value = 1
$students.Cells.Item(5,3).Value = 1
$students.Cells.Item(5,4).Value = "ISR-3"

$students.Columns.Item(1).ColumnWidth = 18.14
$students.Columns.Item(2).ColumnWidth = 9.57
$students.Columns.Item(3).ColumnWidth = 9.86

# ---------------------------------------------------------------------
# 3. Instructors sheet: un-hide MinHours/MaxHours columns, drop the
#    student-crew columns (now living on the "Students" sheet), and
#    bump a handful of required-staff counts.
# ---------------------------------------------------------------------
$instructors.Columns.Item(4).Hidden = $false
$instructors.Columns.Item(5).Hidden = $false

$instructors.Columns.Item(7).ClearContents()
$instructors.Columns.Item(9).ClearContents()

$instructors.Range("C2").Value = 2
$instructors.Range("C3").Value = 2
$instructors.Range("B5").Value = 2
$instructors.Range("C5").Value = 2
$instructors.Range("B6").Value = 2
$instructors.Range("C6").Value = 2
$instructors.Range("C8").Value = 2

# ---------------------------------------------------------------------
# 4. Populate "Sheet3" with the Event / ValueForCompleting table.
# ---------------------------------------------------------------------
$sheet3 = $wb.Worksheets.Item("Sheet3")

$sheet3.Cells.Item(1,1).Value = "Event"
$sheet3.Cells.Item(1,2).Value = "ValueForCompleting"

$events = @(
    @("TR-1", 1),
    @("TR-2", 1),
    @("TR-3", 1),
    @("TR-4", 1),
    @("TR-5", 1),
    @("ISR-1", 2),
    @("ISR-2", 2),
    @("ISR-3", 2),
    @("ISR-4", 2),
    @("ISR-5", 2),
    @("SCAR-1", 3),
    @("SCAR-2", 3),
    @("CAS-1", 4),
    @("CAS-2", 4),
    @("CAS-3", 4),
    @("CAS-4", 4),
    @("CAS-5", 4),
    @("CAS-6", 4)
)

$r = 2
foreach ($row in $events) {
    $sheet3.Cells.Item($r, 1).Value = $row[0]
    $sheet3.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Keep "Instructors" as the selected/active tab, matching the source tab
# state (it was the active sheet before these edits too).
$instructors.Activate()
